$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.07729160785675
$ws.Range("B1").Value = 2.404224395751953
$ws.Range("C1").Value = 6.509651660919189
$ws.Range("D1").Value = 2.229686737060547
$ws.Range("E1").Value = 1.282644629478455
